# Inicializacion de la base de datos
# Adds two more "blank" employee rows (7 and 8) to the sheet, continuing the
# id sequence in column A (6, 7). All the other columns (B..T) are left
# empty, matching the existing placeholder rows already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($col in $cols) {
    $ws.Range($col + "7").Value = ""
    $ws.Range($col + "8").Value = ""
}
